$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume/1h changes) per the diff.
# Rows 30/31 and 48/49 also swap which coin occupies which row.
# The Price column (D) holds plain-looking numeric strings (e.g. '138.30',
# '0.0676', thousand-dot formatted big prices like '56.598.74') that Excel
# must keep as literal TEXT (matching the source file's inline-string cells),
# so NumberFormat is forced to Text ('@') before each Price write.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.598.74'
$ws.Range("E2").Value = '  +2.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.995.34'
$ws.Range("E3").Value = '  +1.69%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '509.87'
$ws.Range("E5").Value = '  +5.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.30'
$ws.Range("E6").Value = '  +6.77%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.434'
$ws.Range("E8").Value = '  +4.89%  '

$ws.Range("E9").Value = '  +7.11%  '

$ws.Range("E10").Value = '  +8.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.355'
$ws.Range("E11").Value = '  +3.00%  '

$ws.Range("E12").Value = '  +2.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.509.41'
$ws.Range("E13").Value = '  +1.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.71'
$ws.Range("E14").Value = '  +5.51%  '

$ws.Range("E15").Value = '  +13.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.626.20'
$ws.Range("E16").Value = '  +2.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.996.23'
$ws.Range("E17").Value = '  +1.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.93'
$ws.Range("E18").Value = '  +6.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.46'
$ws.Range("E19").Value = '  +4.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.82'
$ws.Range("E20").Value = '  +6.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.42'
$ws.Range("E21").Value = '  +5.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").Value = '  -0.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.484'
$ws.Range("E23").Value = '  +5.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.10'
$ws.Range("E24").Value = '  +5.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.170'
$ws.Range("E25").Value = '  +7.50%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0898'
$ws.Range("E27").Value = '  +8.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.63'
$ws.Range("E28").Value = '  +2.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.94'
$ws.Range("E29").Value = '  +7.89%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.80'
$ws.Range("E30").Value = '  +8.03%  '

$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.21'
$ws.Range("E31").Value = '  +5.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.63'
$ws.Range("E32").Value = '  +7.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '154.23'
$ws.Range("E33").Value = '  +4.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.52'
$ws.Range("E34").Value = '  +4.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.66'
$ws.Range("E35").Value = '  +1.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.26'
$ws.Range("E36").Value = '  -0.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0676'
$ws.Range("E37").Value = '  +5.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.48'
$ws.Range("E38").Value = '  +1.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.028.91'
$ws.Range("E39").Value = '  +1.72%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.94'
$ws.Range("E40").Value = '  +3.23%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.271.35'
$ws.Range("E42").Value = '  +7.48%  '

$ws.Range("E43").Value = '  +3.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.68'
$ws.Range("E44").Value = '  +5.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.992'
$ws.Range("E46").Value = '  +1.16%  '

$ws.Range("E47").Value = '  +12.43%  '

$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.84'
$ws.Range("E48").Value = '  +5.86%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0237'
$ws.Range("E49").Value = '  +3.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.15'
$ws.Range("E50").Value = '  +1.98%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0872'
$ws.Range("E51").Value = '  +5.98%  '
